$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total)
$ws.Range("B12").Value = 75
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "72.6/140"
